$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: firstName, lastName, username updates
$ws.Range("A2").Value = "ivab"
$ws.Range("C2").Value = "antons"
$ws.Range("E2").Value = "ivab123"

# Row 3: firstName, lastName, username updates
$ws.Range("A3").Value = "nikod"
$ws.Range("C3").Value = "anilf"
$ws.Range("E3").Value = "nikod321"

# Row 4: firstName, lastName, username updates
$ws.Range("A4").Value = "loris"
$ws.Range("C4").Value = "hoyw"
$ws.Range("E4").Value = "lorisl321"
